# RTM.xlsx update: add two new test cases (5a/5b/6), a "Defects" column,
# and a Defect/Critical Failure summary box, per commit:
#   "Updated.  Found a defect and a critical failure."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- colour constants (BGR ints for OLE Color) ----
$colPink   = 10040319   # FFFF3399 - existing header fill
$colRed    = 255        # FFFF0000 - "Failure" fill
$colDarkRed = 192       # FFC00000 - "Critical Failure" fill

$hLeft   = -4131
$vCenter = -4108

# ---------------------------------------------------------------------
# Row 1 - small Defect / Critical Failure legend box next to the title
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Defect"
$ws.Range("B1").Interior.Color = $colRed

$ws.Range("C1").Value = "Critical Failure"
$ws.Range("C1").Interior.Color = $colDarkRed

# ---------------------------------------------------------------------
# Row 3 - header row, now spanning A:F (was B:F) with a new "Defects" col
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "Test Case"
$ws.Range("B3").Value = "Scenario"
$ws.Range("C3").Value = "Steps"
$ws.Range("D3").Value = "Test Data"
$ws.Range("E3").Value = "Result"
$ws.Range("F3").Value = "Defects"

$headerRow = $ws.Range("A3:F3")
$headerRow.Interior.Color = $colPink
$headerRow.HorizontalAlignment = $hLeft

# ---------------------------------------------------------------------
# Existing test cases 1-4 shift from columns B:F into A:F
# ---------------------------------------------------------------------
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Log In"
$ws.Range("C4").Value = "1) Open Application`n2) Enter Login Credentials"
$ws.Range("D4").Value = "id = kurt`npw = kurt"
$ws.Range("E4").Value = "Pass"
$ws.Range("F4").ClearContents() | Out-Null

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Add Account to Database`nFor Existing Users"
$ws.Range("C5").Value = "1) Enter Account Details`n2) Click ""Add Account To Database"""
$ws.Range("D5").Value = "id = hello2`npwd = 123123"
$ws.Range("E5").Value = "Pass"
$ws.Range("F5").ClearContents() | Out-Null

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "Maintain New Account`nData after logout"
$ws.Range("C6").Value = "1) Log out`n2) Log In`n3) Check for hello2 account"
$ws.Range("D6").ClearContents() | Out-Null
$ws.Range("E6").Value = "Pass"
$ws.Range("F6").ClearContents() | Out-Null

$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Check Account Persistance`nBetween Application`nLaunches"
$ws.Range("C7").Value = "1) Log In`n2) Add Account`n3) Log Out`n4) Log In`n5) Check for New Account"
$ws.Range("D7").Value = "id = hello2`npwd = 123123"
$ws.Range("E7").Value = "Pass"
$ws.Range("F7").ClearContents() | Out-Null

# formatting for rows 4-7: column A/E plain centred, B/C/D wrapped
$plainCols = $ws.Range("A4:A7,E4:E7")
$plainCols.HorizontalAlignment = $hLeft
$plainCols.VerticalAlignment = $vCenter

$wrapCols = $ws.Range("B4:D7")
$wrapCols.HorizontalAlignment = $hLeft
$wrapCols.VerticalAlignment = $vCenter
$wrapCols.WrapText = $true

$ws.Range("D6").HorizontalAlignment = $hLeft
$ws.Range("D6").VerticalAlignment = $vCenter

$ws.Range("F4:F7").HorizontalAlignment = $hLeft
$ws.Range("F4:F7").VerticalAlignment = $vCenter

# ---------------------------------------------------------------------
# Row 8 - test case 5 summary line (spans into 5a/5b below)
# ---------------------------------------------------------------------
$ws.Range("A8").Value = 5
$ws.Range("A8").HorizontalAlignment = $hLeft
$ws.Range("A8").VerticalAlignment = $vCenter

$ws.Range("B8").Value = "Register New User and Add`nAccount Info"
$ws.Range("B8").HorizontalAlignment = $hLeft
$ws.Range("B8").VerticalAlignment = $vCenter
$ws.Range("B8").WrapText = $true

$ws.Range("C8:F8").HorizontalAlignment = $hLeft
$ws.Range("C8:F8").VerticalAlignment = $vCenter

# ---------------------------------------------------------------------
# Row 9 - test case 5a: Register New User
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "5a"
$ws.Range("B9").Value = "Register New User"
$ws.Range("C9").Value = "1) Register New User"
$ws.Range("D9").Value = "id = new`npwd = new"
$ws.Range("E9").Value = "Pass"

$ws.Range("A9,B9,C9,E9,F9").HorizontalAlignment = $hLeft
$ws.Range("A9,B9,C9,E9,F9").VerticalAlignment = $vCenter
$ws.Range("D9").HorizontalAlignment = $hLeft
$ws.Range("D9").VerticalAlignment = $vCenter
$ws.Range("D9").WrapText = $true

# ---------------------------------------------------------------------
# Row 10 - test case 5b: Add Account to Database For Existing Users
#          (highlighted red - a defect was logged against it)
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "5b"
$ws.Range("B10").Value = "Add Account Info For New`nUser"
$ws.Range("C10").Value = "1) Log in as new user without `nany into in database`n2) Add account info"
$ws.Range("D10").Value = "id = 444`npwd = 444"
$ws.Range("E10").Value = "Failure"
$ws.Range("F10").Value = "Defect_01:`nMust Enter account`nInto multuple times`nto save"

$row10red = $ws.Range("A10:F10")
$row10red.Interior.Color = $colRed

$wrapRed = $ws.Range("B10:D10,F10")
$wrapRed.HorizontalAlignment = $hLeft
$wrapRed.VerticalAlignment = $vCenter
$wrapRed.WrapText = $true

$ws.Range("E10").HorizontalAlignment = $hLeft
$ws.Range("E10").VerticalAlignment = $vCenter

# ---------------------------------------------------------------------
# Row 11 - test case 6: Register Multiple Users In One Session
#          (highlighted dark red - a critical failure/crash was logged)
# ---------------------------------------------------------------------
$ws.Range("A11").Value = 6
$ws.Range("B11").Value = "Register Multiple Users In`nOne Sessions"
$ws.Range("C11").Value = "1) Register New User`n2) Log In`n3) Log Out`n4) Register New User"
$ws.Range("D11").Value = "id = Hello`npwd = Hello"
$ws.Range("E11").Value = "Failure"
$ws.Range("F11").Value = "Crash_01:`nProgram crash upon second`naccount creation"

$row11dark = $ws.Range("A11:F11")
$row11dark.Interior.Color = $colDarkRed

$ws.Range("A11").HorizontalAlignment = $hLeft
$ws.Range("A11").VerticalAlignment = $vCenter
$ws.Range("E11").HorizontalAlignment = $hLeft
$ws.Range("E11").VerticalAlignment = $vCenter

$wrapDark = $ws.Range("B11:D11,F11")
$wrapDark.HorizontalAlignment = $hLeft
$wrapDark.VerticalAlignment = $vCenter
$wrapDark.WrapText = $true

# ---------------------------------------------------------------------
# Row heights (auto-fit in the authored workbook; pin explicitly so the
# wrapped text renders the same regardless of host font metrics)
# ---------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 75
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60

# ---------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 26
$ws.Columns.Item(3).ColumnWidth = 33
$ws.Columns.Item(5).ColumnWidth = 13.5
$ws.Columns.Item(6).ColumnWidth = 28.166666666666668

# ---------------------------------------------------------------------
# Selection moves from the old G3 to C1
# ---------------------------------------------------------------------
$ws.Range("C1").Select() | Out-Null
